# Applies Russian translations to GRAPHICS.docx per the commit
# "New translations GRAPHICS.docx (Russian)".
$d = $word.ActiveDocument
$nbsp = [char]0x00A0

function Replace-Text($find, $replace) {
    $result = $d.Content.Find.Execute(
        $find, $true, $false, $false, $false, $false,
        $true, 1, $false, $replace, 2)
    if (-not $result) {
        Write-Output ("NOT FOUND: " + $find)
    }
}

Replace-Text "GRAPHICS RESOURCES" "БИБЛИОТЕКА ГРАФИКИ"

Replace-Text "These graphics can be used for display in online services, retail establishments, press or other promotional purposes" "Вы можете использовать представленную ниже графику для онлайн сервисов, торговых заведений, прессы и т.д."

$logosFind = "LOGOS" + $nbsp
$logosReplace = "ЛОГОТИПЫ" + $nbsp
Replace-Text $logosFind $logosReplace

$iconsFind = $nbsp + "ICONS"
$iconsReplace = $nbsp + "ИКОНКИ"
Replace-Text $iconsFind $iconsReplace

Replace-Text "MERCHANT BUTTONS" "ДЛЯ МЕРЧАНТОВ"

Replace-Text "SMARTCASH ROADMAP FEATURES" "ДОРОЖНАЯ КАРТА: ФУНКЦИИ"

Replace-Text "SOCIAL MEDIA GRAPHICS" "ГРАФИКА ДЛЯ СОЦИАЛЬНЫХ СЕТЕЙ"

Replace-Text " GRAPHIC IDENTITY" " ХАРАКТЕРИСТИКИ"

Replace-Text "Graphic Identity Guidelines" "Рекомендуемые параметры для графических материалов"

$guidelinesFind = "These are recommended usage guidelines for maintaining a consistent design aesthetic for the SmartCash brand." + $nbsp + "A strong and consistent visual identity of our logo will help keep a consistent look, recognition and familiarity now and in the future. Standardization of colours will go a long way to enforce a reliable and positive impression to our identity in the blockchain space."
$guidelinesReplace = "Эти рекомендации созданы для поддержания дизайна бренда SmartCash.Точная и последовательная визуальная идентификация нашего логотипа поможет нам оставаться узнаваемыми в криптовалютном пространстве сейчас и в будущем. Standardization of colours will go a long way to enforce a reliable and positive impression to our identity in the blockchain space."
Replace-Text $guidelinesFind $guidelinesReplace

Write-Output "Done."
